$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells D1:F1 - copy formatting (style) from existing header A1
$ws.Range("A1").Copy()
$ws.Range("D1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("D1").Value = "产品名称2"
$ws.Range("E1").Value = "平台名称2"
$ws.Range("F1").Value = "备注2"

$data = @(
  @("瑞e保", "上海中振交通装备有限公司", "上海中振"),
  @("瑞e保", "上海电气研砼建筑科技集团有限公司", "上海电气研砼"),
  @("瑞e仓", "上海中钢供应链有限公司", "中钢"),
  @("瑞e惠", "云度新能源汽车有限公司", "云度汽车"),
  @("瑞e险", "熠瑞合（上海）供应链管理有限公司", "壹链盟"),
  @("瑞e订", "上海大展懿杰建设科技有限公司", "大展"),
  @("瑞e惠", "平安普惠融资担保公司", "平安普惠"),
  @("瑞e惠", "恒洁卫浴集团有限公司", "恒洁卫浴"),
  @("瑞e惠", "磁金融", "磁金融"),
  @("瑞e保", "江西联创光电科技股份有限公司", "联创光电"),
  @("瑞e保", "远东智慧能源股份有限公司", "远东"),
  @("瑞e订", "钰翔环保供应链管理（上海）有限公司", "钰翔"),
  @("瑞e保", "银江技术股份有限公司", "银江"),
  @("瑞e惠", "锱云（上海）物联网科技有限公司", "锱云"),
  @("瑞e保", "龙元供应链管理有限公司", "龙元建设")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $data[$i][0]
    $ws.Cells.Item($row, 5).Value = $data[$i][1]
    $ws.Cells.Item($row, 6).Value = $data[$i][2]
}
